$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 10 ("Parametros de tsconfig.json"),
# pushing it (and the "sourceMap" row after it) down by one. This leaves
# row 9 free for a new outFile usage example and row 10 empty, matching
# the original blank-row separator.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "tsc outFile archivoSalida archivo1 archivo2 archivo3"

# The "sourceMap" row (now row 12) gains a "Tipo boolean" column
$ws.Range("C12").Value = "Tipo boolean"

# Row 13: removeComments
$ws.Range("A13").Value = "removeComments"
$ws.Range("B13").Value = "Permite indicar si el js se genera con comentarios"
$ws.Range("C13").Value = "Tipo boolean"

# Row 14: include
$ws.Range("A14").Value = "include"
$ws.Range("B14").Value = "Permite indicar que directorios incluir en la compilación, no se define dentro del compileroptions"
$ws.Range("C14").Value = "Arreglo Strings"

# Row 15: exclude
$ws.Range("A15").Value = "exclude"
$ws.Range("B15").Value = "Permite indicar que directorios excluir en la compilación, no se define dentro del compileroptions"
$ws.Range("C15").Value = "Arreglo Strings"

# Row 16: outFile
$ws.Range("A16").Value = "outFile"
$ws.Range("B16").Value = "en el tsconfig.json recibe la dirección y nombre archivo salida de un compilado de todos los ts files, además se requiere compilar a mano con tsc."
$ws.Range("C16").Value = "String"

# Row 17: file (with hyperlink in D17)
$ws.Range("A17").Value = "file"
$ws.Range("B17").Value = "El parametro file permite indicar los archivos a compilar y esto tambien ayuda a determinar el orden de generación de los archivos en el outfile."
$ws.Range("C17").Value = "Arreglo Strings"
$ws.Range("D17").Value = "http://www.typescriptlang.org/docs/handbook/tsconfig-json.html"
$ws.Hyperlinks.Add($ws.Range("D17"), "http://www.typescriptlang.org/docs/handbook/tsconfig-json.html")

# Widen column B to fit the longer descriptions
$ws.Columns.Item(2).ColumnWidth = 90.140625

# Update selection to mirror the saved state from the source workbook
$ws.Range("B17").Select()
